$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 352.875
$ws.Range("I6").Value = 324.95
$ws.Range("K6").Value = 974.8499999999999
$ws.Range("M6").Value = -862.8499999999999
# Row 41
$ws.Range("H41").Value = 2205
$ws.Range("I41").Value = 2571.7778
$ws.Range("J41").Value = 1654.8334
$ws.Range("K41").Value = 2571.7778
$ws.Range("L41").Value = 1654.8334
$ws.Range("M41").Value = -2131.7778
$ws.Range("N41").Value = -2534.8334
# Row 115
$ws.Range("H115").Value = 6665.3335
$ws.Range("J115").Value = 10000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -33134
# Row 131
$ws.Range("H131").Value = 608.8125
$ws.Range("I131").Value = 608.8125
$ws.Range("K131").Value = 1826.4375
$ws.Range("M131").Value = 3213.5625
# Row 132
$ws.Range("H132").Value = 1103.4642
$ws.Range("I132").Value = 1122.8518
$ws.Range("K132").Value = 3368.5554
$ws.Range("M132").Value = -838.5553999999997
# Row 138
$ws.Range("H138").Value = 2104.389
$ws.Range("J138").Value = 3133.111
$ws.Range("L138").Value = 9399.332999999999
$ws.Range("N138").Value = -19679.333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 4363.1665
$ws.Range("I26").Value = 4500
$ws.Range("J26").Value = 4294.75
$ws.Range("K26").Value = 4500
$ws.Range("L26").Value = 4294.75
$ws.Range("M26").Value = -4170
$ws.Range("N26").Value = -4954.75
# Row 33
$ws.Range("H33").Value = 6791.3335
$ws.Range("I33").Value = 6791.3335
$ws.Range("K33").Value = 6791.3335
$ws.Range("M33").Value = -6462.3335
# Row 97
$ws.Range("H97").Value = 1916.1666
$ws.Range("I97").Value = 1518.2106
$ws.Range("J97").Value = 3428.4
$ws.Range("K97").Value = 1518.2106
$ws.Range("L97").Value = 3428.4
$ws.Range("M97").Value = -1022.2106
$ws.Range("N97").Value = -4420.4
# Row 132
$ws.Range("H132").Value = 1115.1765
$ws.Range("I132").Value = 1028.625
$ws.Range("K132").Value = 3085.875
$ws.Range("M132").Value = -555.875

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 7944.0713
$ws.Range("J20").Value = 2581.75
$ws.Range("L20").Value = 2581.75
$ws.Range("N20").Value = -3075.75
# Row 86
$ws.Range("H86").Value = 8588.046
$ws.Range("I86").Value = 2715
$ws.Range("J86").Value = 11944.071
$ws.Range("K86").Value = 2715
$ws.Range("L86").Value = 11944.071
$ws.Range("M86").Value = -1592
$ws.Range("N86").Value = -14190.071
# Row 89
$ws.Range("H89").Value = 8588.046
$ws.Range("I89").Value = 2715
$ws.Range("J89").Value = 11944.071
$ws.Range("K89").Value = 13575
$ws.Range("L89").Value = 59720.355
$ws.Range("M89").Value = -7959
$ws.Range("N89").Value = -70952.355
# Row 134
$ws.Range("H134").Value = 8125.727
$ws.Range("I134").Value = 8125.727
$ws.Range("K134").Value = 24377.181
$ws.Range("M134").Value = -21842.181

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 13000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 31
$ws.Range("H31").Value = 1971.1945
$ws.Range("I31").Value = 2017.7778
$ws.Range("J31").Value = 1924.6111
$ws.Range("K31").Value = 2017.7778
$ws.Range("L31").Value = 1924.6111
$ws.Range("M31").Value = -1722.7778
$ws.Range("N31").Value = -2514.6111
# Row 34
$ws.Range("H34").Value = 1971.1945
$ws.Range("I34").Value = 2017.7778
$ws.Range("J34").Value = 1924.6111
$ws.Range("K34").Value = 2017.7778
$ws.Range("L34").Value = 1924.6111
$ws.Range("M34").Value = -1815.7778
$ws.Range("N34").Value = -2328.6111
# Row 38
$ws.Range("H38").Value = 23495
$ws.Range("I38").Value = 16990
$ws.Range("J38").Value = 30000
$ws.Range("K38").Value = 16990
$ws.Range("L38").Value = 30000
$ws.Range("M38").Value = -16613
$ws.Range("N38").Value = -30754
# Row 46
$ws.Range("H46").Value = 23495
$ws.Range("I46").Value = 16990
$ws.Range("J46").Value = 30000
$ws.Range("K46").Value = 16990
$ws.Range("L46").Value = 30000
$ws.Range("M46").Value = -16779
$ws.Range("N46").Value = -30422
# Row 99
$ws.Range("H99").Value = 2109.5
$ws.Range("I99").Value = 1974.6666
$ws.Range("K99").Value = 1974.6666
$ws.Range("M99").Value = -476.6666
# Row 126
$ws.Range("H126").Value = 2109.5
$ws.Range("I126").Value = 1974.6666
$ws.Range("K126").Value = 5923.9998
$ws.Range("M126").Value = -3453.9998
# Row 134
$ws.Range("H134").Value = 2231
$ws.Range("I134").Value = 2177.4443
$ws.Range("K134").Value = 6532.3329
$ws.Range("M134").Value = -3997.3329

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 8481.75
$ws.Range("J18").Value = 2976
$ws.Range("L18").Value = 8928
$ws.Range("N18").Value = -9266
# Row 39
$ws.Range("H39").Value = 5042.857
$ws.Range("I39").Value = 100
$ws.Range("J39").Value = 5866.6665
$ws.Range("K39").Value = 300
$ws.Range("L39").Value = 17599.9995
$ws.Range("M39").Value = -6
$ws.Range("N39").Value = -18187.9995
# Row 56
$ws.Range("H56").Value = 3666.6667
$ws.Range("I56").Value = 3666.6667
$ws.Range("K56").Value = 3666.6667
$ws.Range("M56").Value = -3136.6667
# Row 113
$ws.Range("H113").Value = 583.4167
$ws.Range("I113").Value = 551
$ws.Range("J113").Value = 589.9
$ws.Range("K113").Value = 1653
$ws.Range("L113").Value = 1769.7
$ws.Range("M113").Value = 517
$ws.Range("N113").Value = -6109.7
# Row 122
$ws.Range("H122").Value = 775.82355
$ws.Range("I122").Value = 506.92307
$ws.Range("J122").Value = 1649.75
$ws.Range("K122").Value = 4562.30763
$ws.Range("L122").Value = 14847.75
$ws.Range("M122").Value = -2112.30763
$ws.Range("N122").Value = -19747.75
# Row 134
$ws.Range("H134").Value = 21700
$ws.Range("I134").Value = 21700
$ws.Range("K134").Value = 65100
$ws.Range("M134").Value = -60030
# Row 136
$ws.Range("H136").Value = 4689.75
$ws.Range("J136").Value = 6975
$ws.Range("L136").Value = 20925
$ws.Range("N136").Value = -31125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 6993.3335
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888
# Row 8
$ws.Range("H8").Value = 6993.3335
$ws.Range("I8").Value = 3000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2861
# Row 19
$ws.Range("H19").Value = 19999
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
# Row 80
$ws.Range("H80").Value = 2292.1428
$ws.Range("I80").Value = 2295
$ws.Range("J80").Value = 2291
$ws.Range("K80").Value = 2295
$ws.Range("L80").Value = 2291
$ws.Range("M80").Value = -1297
$ws.Range("N80").Value = -4287
# Row 83
$ws.Range("H83").Value = 2292.1428
$ws.Range("I83").Value = 2295
$ws.Range("J83").Value = 2291
$ws.Range("K83").Value = 11475
$ws.Range("L83").Value = 11455
$ws.Range("M83").Value = -6483
$ws.Range("N83").Value = -21439
# Row 97
$ws.Range("H97").Value = 744.0769
$ws.Range("I97").Value = 519.6667
$ws.Range("K97").Value = 519.6667
$ws.Range("M97").Value = -23.66669999999999
# Row 113
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3500
$ws.Range("N113").Value = -7840
$ws.Range("M113").ClearContents()
# Row 132
$ws.Range("H132").Value = 1719.3
$ws.Range("I132").Value = 1702.238
$ws.Range("J132").Value = 1759.1111
$ws.Range("K132").Value = 5106.714
$ws.Range("L132").Value = 5277.3333
$ws.Range("M132").Value = -2576.714
$ws.Range("N132").Value = -10337.3333
# Row 136
$ws.Range("H136").Value = 31512.834
$ws.Range("J136").Value = 31512.834
$ws.Range("L136").Value = 94538.50199999999
$ws.Range("N136").Value = -99638.50199999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 11
$ws.Range("H11").Value = 9994.5
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
# Row 20
$ws.Range("H20").Value = 8667.166999999999
$ws.Range("J20").Value = 8667.166999999999
$ws.Range("L20").Value = 8667.166999999999
$ws.Range("N20").Value = -9119.166999999999
# Row 22
$ws.Range("H22").Value = 1456.4667
$ws.Range("I22").Value = 1542.0769
$ws.Range("K22").Value = 1542.0769
$ws.Range("M22").Value = -1247.0769
# Row 27
$ws.Range("H27").Value = 1456.4667
$ws.Range("I27").Value = 1542.0769
$ws.Range("K27").Value = 1542.0769
$ws.Range("M27").Value = -1435.0769
# Row 34
$ws.Range("H34").Value = 24124.75
$ws.Range("I34").Value = 23499
$ws.Range("K34").Value = 23499
$ws.Range("M34").Value = -23327
# Row 55
$ws.Range("H55").Value = 117.73333
$ws.Range("I55").Value = 94
$ws.Range("K55").Value = 94
$ws.Range("M55").Value = 79
# Row 122
$ws.Range("H122").Value = 3099.1875
$ws.Range("I122").Value = 2806.5454
$ws.Range("K122").Value = 8419.636200000001
$ws.Range("M122").Value = -5969.636200000001
# Row 132
$ws.Range("H132").Value = 7306.7
$ws.Range("I132").Value = 6083.4287
$ws.Range("J132").Value = 10161
$ws.Range("K132").Value = 18250.2861
$ws.Range("L132").Value = 30483
$ws.Range("M132").Value = -15720.2861
$ws.Range("N132").Value = -35543

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 3000
$ws.Range("I8").Value = 3000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2860
# Row 113
$ws.Range("H113").Value = 1663.0555
$ws.Range("I113").Value = 540.1818
$ws.Range("K113").Value = 1620.5454
$ws.Range("M113").Value = 549.4546
# Row 122
$ws.Range("H122").Value = 3709.1538
$ws.Range("I122").Value = 3330.652
$ws.Range("K122").Value = 9991.956
$ws.Range("M122").Value = -7541.956
